$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide rows 3-32 and 34-44 (row 33 stays visible)
$ws.Range("A3:A32").EntireRow.Hidden = $true
$ws.Range("A34:A44").EntireRow.Hidden = $true

# Update B45 value from 41 to 43
$ws.Range("B45").Value = 43

# Zoom to 175% and move the selection to B46
$excel.ActiveWindow.Zoom = 175
$ws.Range("B46").Select() | Out-Null
